$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update the "last updated" timestamp in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 11:52"

# --- 2) Insert a new row for "Malasia" before the current Pakistan row (row 33), ---
#        then delete the old Malasia row (now shifted down to row 36) so the data ---
#        table keeps country Malasia but in its new sorted position. ---
$ws.Rows.Item(33).Insert()
$ws.Range("A33").Value = "Malasia"
$ws.Range("B33").Value = 3963
$ws.Range("C33").Value = 170
$ws.Range("D33").Value = 1321
$ws.Range("E33").Value = 2579
$ws.Range("F33").Value = 92
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = 63

# Old Malasia row has shifted from 35 -> 36 because of the insert above.
$ws.Rows.Item(36).Delete()

# --- 3) Insert a new row for "Malta" before the current Ghana row (row 99), ---
#        then delete the old Malta row (now shifted down to row 106). ---
$ws.Rows.Item(99).Insert()
$ws.Range("A99").Value = "Malta"
$ws.Range("B99").Value = 293
$ws.Range("C99").Value = 52
$ws.Range("D99").Value = 5
$ws.Range("E99").Value = 288
$ws.Range("F99").Value = 3
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0

# Old Malta row has shifted from 105 -> 106 because of the insert above.
$ws.Rows.Item(106).Delete()

# --- 4) Direct numeric refreshes for countries whose stats changed but whose ---
#        table position did not move. ---

# España (row 5)
$ws.Range("B5").Value = 140510
$ws.Range("C5").Value = 3835
$ws.Range("D5").Value = 43208
$ws.Range("E5").Value = 83504
$ws.Range("F5").Value = 7069
$ws.Range("G5").Value = 457
$ws.Range("H5").Value = 13798

# Suiza (row 14)
$ws.Range("B14").Value = 21855
$ws.Range("C14").Value = 198
$ws.Range("E14").Value = 13015
$ws.Range("G14").Value = 19
$ws.Range("H14").Value = 784

# Brasil (row 18)
$ws.Range("B18").Value = 12240
$ws.Range("C18").Value = 57
$ws.Range("E18").Value = 11547

# India (row 28)
$ws.Range("E28").Value = 4339
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 137

# Hong Kong (row 66)
$ws.Range("B66").Value = 936
$ws.Range("C66").Value = 21
$ws.Range("D66").Value = 236
$ws.Range("E66").Value = 696

# Bosnia y Herzegovina (row 72)
$ws.Range("E72").Value = 640
$ws.Range("G72").Value = 3
$ws.Range("H72").Value = 32

# Libano (row 82)
$ws.Range("B82").Value = 548
$ws.Range("C82").Value = 7
$ws.Range("E82").Value = 469

# Sri Lanka (row 113)
$ws.Range("D113").Value = 42
$ws.Range("E113").Value = 132
